# Insert a new data row (weekly Coco price record) above the existing row 32,
# shifting the subsequent rows down by one (row 32 -> 33, ..., row 43 -> 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = 9
$ws.Range("B32").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 45258
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108007
$ws.Range("J32").Value = "Coco"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 150
$ws.Range("N32").Value = 28000
$ws.Range("O32").Value = 28000
$ws.Range("P32").Value = 28000
$ws.Range("Q32").Value = "$/malla 20 unidades"
$ws.Range("R32").Value = "Perú"
$ws.Range("S32").Value = 1400
$ws.Range("T32").Value = 20
